$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: the "Тема" bullet ends with a run containing only ";" -
# change that run's text to ":" (keeping the run / its rPr intact,
# i.e. only the character changes, not the run structure).
# ------------------------------------------------------------------
$topicHit = $d.Content
$topicHit.Find.Execute("Тема;") | Out-Null
$semicolon = $d.Range($topicHit.End - 1, $topicHit.End)
$semicolon.Text = ":"

# ------------------------------------------------------------------
# Edit 2: "Шаблоны" is currently split across two runs - "Шаблон"
# and "ы". Merge them into a single run containing "Шаблоны".
# ------------------------------------------------------------------
$wholeWord = $d.Content
$wholeWord.Find.Execute("Шаблоны") | Out-Null
$wholeEnd = $wholeWord.End

$firstPart = $d.Content
$firstPart.Find.Execute("Шаблон") | Out-Null
$firstPartEnd = $firstPart.End

$tail = $d.Range($firstPartEnd, $wholeEnd)
$tailText = $tail.Text
$tail.Text = ""

$insertPoint = $d.Range($firstPartEnd, $firstPartEnd)
$insertPoint.InsertAfter($tailText)
